$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Resize the table: overall preferred width + the five grid columns.
#    COM widths are expressed in points; OOXML dxa/pct values are twentieths
#    of that, so divide the target OOXML values by 20.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)
$t.PreferredWidth = 4291 / 20               # w:tblW w:w="4291" w:type="pct"

$t.Columns.Item(1).Width = 2649 / 20        # gridCol 1
$t.Columns.Item(2).Width = 1840 / 20        # gridCol 2
$t.Columns.Item(3).Width = 1271 / 20        # gridCol 3
$t.Columns.Item(4).Width = 1264 / 20        # gridCol 4
$t.Columns.Item(5).Width = 756 / 20         # gridCol 5

# ---------------------------------------------------------------------------
# 2. Header-row label swaps. "N" / "Y" are single letters that also occur as
#    the tail of unrelated strings (e.g. "donor_sex_mismatch: Y"), so match
#    only the very first whole-word hit starting from the top of the doc.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("N", $true, $true, $false, $false, $false, $true, 1, $false,
                         "No/Mild/Mod. PGD (n = 54)", 2)

$d.Content.Find.Execute("Y", $true, $true, $false, $false, $false, $true, 1, $false,
                         "Severe PGD (n = 8)", 1)

# ---------------------------------------------------------------------------
# 3. Re-computed statistics / relabelled variable name.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("0.17", $true, $false, $false, $false, $false, $true, 1, $false,
                         "0.165", 2)

$d.Content.Find.Execute(">0.99", $true, $false, $false, $false, $false, $true, 1, $false,
                         ">0.999", 2)

$d.Content.Find.Execute("0.33", $true, $false, $false, $false, $false, $true, 1, $false,
                         "0.328", 2)

$d.Content.Find.Execute("donor_PHM", $true, $true, $false, $false, $false, $true, 1, $false,
                         "donor_PHM_calc", 2)

$d.Content.Find.Execute("0.79", $true, $false, $false, $false, $false, $true, 1, $false,
                         "0.793", 2)

$d.Content.Find.Execute("0.69", $true, $false, $false, $false, $false, $true, 1, $false,
                         "0.692", 2)

$d.Content.Find.Execute("0.41", $true, $false, $false, $false, $false, $true, 1, $false,
                         "0.414", 2)

# Six identical "0.76" cells all become "0.764" -- replace all occurrences.
$d.Content.Find.Execute("0.76", $true, $false, $false, $false, $false, $true, 1, $false,
                         "0.764", 2)
